$wb = $excel.ActiveWorkbook

$colorBlue = 16770508   # CCE5FF - Equipo/Brawler columns A-C
$colorPink = 13421812   # F4CCCC - Equipo/Brawler columns D-F

function Set-ScrimCell($ws, $addr, $value, $fillColor, $bold) {
    $c = $ws.Range($addr)
    $c.Value = $value
    $c.Borders.LineStyle = 1
    if ($fillColor -ne $null) {
        $c.Interior.Color = $fillColor
    }
    if ($bold) {
        $c.Font.Bold = $true
    }
}

function Add-ScrimRow($ws, $row, $vals) {
    # $vals is an array of 14 values for columns A..N, in order
    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N")
    for ($i = 0; $i -lt 6; $i++) {
        $fill = $colorBlue
        if ($i -ge 3) { $fill = $colorPink }
        Set-ScrimCell $ws "$($cols[$i])$row" $vals[$i] $fill $false
    }
    $winner = $vals[6]
    $winnerBoldFill = $colorPink
    if ($winner -eq "Equipo 1") { $winnerBoldFill = $colorBlue }
    Set-ScrimCell $ws "G$row" $winner $winnerBoldFill $true
    for ($i = 7; $i -lt 14; $i++) {
        Set-ScrimCell $ws "$($cols[$i])$row" $vals[$i] $null $false
    }
}

$ws = $wb.Worksheets.Item("Triple Dribble")
Add-ScrimRow $ws 77 @("MORTIS", "KAZE", "CROW", "LUMI", "BARLEY", "CORDELIUS", "Equipo 2", "IC|Mebius", "IC|Nob?", "IC|RamaZR", "FUT|GeRo", "FUT|Nowy297", "FUT|MeOw", "20250725T174703.000Z")
Add-ScrimRow $ws 78 @("MORTIS", "KAZE", "CROW", "LUMI", "BARLEY", "CORDELIUS", "Equipo 2", "IC|Mebius", "IC|Nob?", "IC|RamaZR", "FUT|GeRo", "FUT|Nowy297", "FUT|MeOw", "20250725T174520.000Z")
Add-ScrimRow $ws 79 @("MORTIS", "KAZE", "CROW", "LUMI", "BARLEY", "CORDELIUS", "Equipo 1", "IC|Mebius", "IC|Nob?", "IC|RamaZR", "FUT|GeRo", "FUT|Nowy297", "FUT|MeOw", "20250725T174143.000Z")
Add-ScrimRow $ws 80 @("DRACO", "FINX", "SHADE", "KENJI", "CORDELIUS", "LARRY & LAWRIE", "Equipo 1", "LOUD|FireCrow", "LOUD|Edinho", "LOUD|KaioDog", "Bielz", "GO|Yichy❦", "Tilo🍥", "20250725T175130.000Z")
Add-ScrimRow $ws 81 @("DRACO", "FINX", "SHADE", "KENJI", "CORDELIUS", "LARRY & LAWRIE", "Equipo 2", "LOUD|FireCrow", "LOUD|Edinho", "LOUD|KaioDog", "Bielz", "GO|Yichy❦", "Tilo🍥", "20250725T174908.000Z")
Add-ScrimRow $ws 82 @("DRACO", "FINX", "SHADE", "KENJI", "CORDELIUS", "LARRY & LAWRIE", "Equipo 1", "LOUD|FireCrow", "LOUD|Edinho", "LOUD|KaioDog", "Bielz", "GO|Yichy❦", "Tilo🍥", "20250725T174655.000Z")
Add-ScrimRow $ws 83 @("JACKY", "MEEPLE", "FRANK", "R-T", "SHADE", "CORDELIUS", "Equipo 2", "LOUD|Edinho", "LOUD|KaioDog", "LOUD|FireCrow", "Bielz", "Tilo🍥", "GO|Yichy❦", "20250725T174051.000Z")
Add-ScrimRow $ws 84 @("JACKY", "MEEPLE", "FRANK", "R-T", "SHADE", "CORDELIUS", "Equipo 1", "LOUD|Edinho", "LOUD|KaioDog", "LOUD|FireCrow", "Bielz", "Tilo🍥", "GO|Yichy❦", "20250725T173759.000Z")

$ws = $wb.Worksheets.Item("Ring of Fire")
Add-ScrimRow $ws 66 @("LOU", "JAE-YONG", "DRACO", "STU", "BEA", "GRIFF", "Equipo 2", "IC|Mebius", "IC|Nob?", "IC|RamaZR", "FUT|GeRo", "FUT|Nowy297", "FUT|MeOw", "20250725T175704.000Z")
Add-ScrimRow $ws 67 @("LOU", "JAE-YONG", "DRACO", "STU", "BEA", "GRIFF", "Equipo 2", "IC|Mebius", "IC|Nob?", "IC|RamaZR", "FUT|GeRo", "FUT|Nowy297", "FUT|MeOw", "20250725T175410.000Z")
